# Add a new "Bruker" scheme row into the averaged-intensities table.
# This inserts a full row at sheet row 4 (shifting the existing schemes
# ND Single..RotRing...OmegaMax-30 down by one row) and fills the new
# row with the Bruker scheme's averaged-intensity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 4; this shifts rows 4:48 down to 5:49
# and extends the sheet dimension to A1:S49 automatically.
$ws.Rows.Item(4).Insert()

# Copy the formatting (bold/border style) from the row above (row 3)
# onto the newly inserted row 4 so column A keeps its header-style border.
$ws.Range("A3:S3").Copy()
$ws.Range("A4:S4").PasteSpecial(-4122)

# Sequence number + scheme label for the new row.
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Bruker"

# Averaged-intensity values for the Bruker scheme.
$ws.Range("C4").Value = 1.012637720645599
$ws.Range("D4").Value = 0.9854906548340724
$ws.Range("E4").Value = 0.9923442785986791
$ws.Range("F4").Value = 1.012637720645599
$ws.Range("G4").Value = 0.9862971980807925
$ws.Range("H4").Value = 0.9845002611835262
$ws.Range("I4").Value = 0.9922589580872768
$ws.Range("J4").Value = 0.9854906548340724
$ws.Range("K4").Value = 1.012637720645599
$ws.Range("L4").Value = 0.9889174667163758
$ws.Range("M4").Value = 0.9889174667163758
$ws.Range("N4").Value = 0.9880440438378481
$ws.Range("O4").Value = 0.9968242180261168
$ws.Range("P4").Value = 0.9968242180261168
$ws.Range("Q4").Value = 1.000777593680987
$ws.Range("R4").Value = 1.000777593680987
$ws.Range("S4").Value = 0.9922548452383243
